# Generate Report for Handback
#
# The file "5b4e3832-64ea-4e4a-b420-0dd2e6015c8e.md" has been handed back
# (its localized translations are in sync with en-US again). Update the
# status everywhere it is tracked, and record the handback target file +
# handback datetime on the per-locale sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: flip the status for both locale columns ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = "Handed back: in sync with en-US"
$overview.Range("C2").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet: status + newly-populated handback columns ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Handed back: in sync with en-US"
$zhcn.Hyperlinks.Add(
    $zhcn.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/d681a4b2a42a176f6298fa088f65ca86405ab641/e2e/5b4e3832-64ea-4e4a-b420-0dd2e6015c8e.md",
    "",
    "",
    "5b4e3832-64ea-4e4a-b420-0dd2e6015c8e.md"
)
$zhcn.Hyperlinks.Add(
    $zhcn.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/87edda102e019bcda82fcb4d76cbedf618efd037/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/5b4e3832-64ea-4e4a-b420-0dd2e6015c8e.6712b26e57c0f76a1b5e98528851d1c8b9b2a095.zh-cn.xlf",
    "",
    "",
    "5b4e3832-64ea-4e4a-b420-0dd2e6015c8e.6712b26e57c0f76a1b5e98528851d1c8b9b2a095.zh-cn.xlf"
)
$zhcn.Range("H2").Value = "2016-03-22 02:32:53"

# --- de-de sheet: status + newly-populated handback columns ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "Handed back: in sync with en-US"
$dede.Hyperlinks.Add(
    $dede.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/d681a4b2a42a176f6298fa088f65ca86405ab641/e2e/5b4e3832-64ea-4e4a-b420-0dd2e6015c8e.md",
    "",
    "",
    "5b4e3832-64ea-4e4a-b420-0dd2e6015c8e.md"
)
$dede.Hyperlinks.Add(
    $dede.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9edb4207171171b77a71716c9b22369a95c6849e/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/5b4e3832-64ea-4e4a-b420-0dd2e6015c8e.6712b26e57c0f76a1b5e98528851d1c8b9b2a095.de-de.xlf",
    "",
    "",
    "5b4e3832-64ea-4e4a-b420-0dd2e6015c8e.6712b26e57c0f76a1b5e98528851d1c8b9b2a095.de-de.xlf"
)
$dede.Range("H2").Value = "2016-03-22 02:32:59"
